$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge "Foi desenvolvido uma aplicação web " + "chamada "Busca
# CEP"" into a single run, without bleeding into the following
# " no Sublime " run (which must stay separate). We place a temporary
# barrier bookmark right before " no Sublime " so the engine's
# run-coalescing pass stops there, then remove the barrier again.
# ---------------------------------------------------------------------
$t = $d.Content.Text
$idxSublime = $t.IndexOf("no Sublime")
$barrier1 = $d.Range($idxSublime, $idxSublime)
$d.Bookmarks.Add("ZZ_BARRIER1", $barrier1)

$t = $d.Content.Text
$idxFoi = $t.IndexOf("Foi desenvolvido uma aplica")
$lenFoi = "Foi desenvolvido uma aplicação web chamada “Busca CEP”".Length
$rFoi = $d.Range($idxFoi, $idxFoi + $lenFoi)
$rFoi.Find.Execute("Foi desenvolvido uma aplicação web chamada “Busca CEP”", $true, $false, $false, $false, $false, $true, 1, $false, "Foi desenvolvido uma aplicação web chamada “Busca CEP”", 2)

$d.Bookmarks("ZZ_BARRIER1").Delete()

# ---------------------------------------------------------------------
# Change 2: remove the "_GoBack" bookmark and merge the ", " run with the
# "onde foi codificado na linguagem HTML5, CSS3, " run. The existing
# "_GoBack" bookmark (still in place at this point) acts as the left
# barrier, and the following <w:proofErr> (around "JavaScript") acts as
# the right barrier, so this merge cannot bleed into " 3 IDE" or
# "JavaScript".
# ---------------------------------------------------------------------
$t = $d.Content.Text
$idxComma = $t.IndexOf(", onde foi codificado na linguagem HTML5, CSS3, ")
$lenComma = ", onde foi codificado na linguagem HTML5, CSS3, ".Length
$rComma = $d.Range($idxComma, $idxComma + $lenComma)
$rComma.Find.Execute(", onde foi codificado na linguagem HTML5, CSS3, ", $true, $false, $false, $false, $false, $true, 1, $false, ", onde foi codificado na linguagem HTML5, CSS3, ", 2)

# ---------------------------------------------------------------------
# Change 3: merge "na linguagem C# no Visual Studio 2019 " + "IDE" + ", "
# into a single run, without touching " desenvolvido um console " before
# it or "para comparar " after it. Temporary barrier bookmarks on both
# sides keep the coalescing pass contained.
# ---------------------------------------------------------------------
$t = $d.Content.Text
$idxNa = $t.IndexOf("na linguagem C# no Visual Studio 2019 ")
$barrier2 = $d.Range($idxNa, $idxNa)
$d.Bookmarks.Add("ZZ_BARRIER2", $barrier2)

$t = $d.Content.Text
$idxPara = $t.IndexOf("para comparar ")
$barrier3 = $d.Range($idxPara, $idxPara)
$d.Bookmarks.Add("ZZ_BARRIER3", $barrier3)

$t = $d.Content.Text
$idxNa = $t.IndexOf("na linguagem C# no Visual Studio 2019 IDE, ")
$lenNa = "na linguagem C# no Visual Studio 2019 IDE, ".Length
$rNa = $d.Range($idxNa, $idxNa + $lenNa)
$rNa.Find.Execute("na linguagem C# no Visual Studio 2019 IDE, ", $true, $false, $false, $false, $false, $true, 1, $false, "na linguagem C# no Visual Studio 2019 IDE, ", 2)

$d.Bookmarks("ZZ_BARRIER2").Delete()
$d.Bookmarks("ZZ_BARRIER3").Delete()

# ---------------------------------------------------------------------
# Change 4: append two new paragraphs at the end of the document — an
# empty one, then one holding the "Obs. ..." note — and relocate the
# "_GoBack" bookmark to sit right after the new note's text (this also
# removes it from its old location, since bookmark names are unique).
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter("Obs. No console ao entrar com os dados, pode ser em caixa alta ou baixa, entretanto deve possuir a acentuação.")

$lastPara = $d.Paragraphs.Last
$finalMark = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $finalMark)
